$wb = $excel.ActiveWorkbook

# Rows in the per-file detail sheets (handoff/handback rows) affected by
# this "generate report for handoff" run.
$rows = @(7, 8, 9, 10, 11, 13)

# --- Overview sheet: "Latest HO Xliff Generate Date" (column G) ---------
# All six rows shared the same handoff timestamp string; bump it forward.
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-08-17 12:20:44"
}

# --- de-de sheet: "Latest Handoff Datetime" (column H) + Priority -------
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Range("H$r").Value = "2016-08-17 12:20:44"
    $wsDeDe.Range("E$r").Value = "ht"
}

# --- zh-cn sheet: "Latest Handoff Datetime" (column H) + Priority -------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Range("H$r").Value = "2016-08-17 12:20:38"
    $wsZhCn.Range("E$r").Value = "ht"
}

